$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D (copy the header style from C1 so D1 matches
# the bold/centered look of the other header cells, then set its text)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "function"

# Data rows 2-8: centers/sigmas become text strings, add function column = "lorentzian"
$data = @(
    @("584.2089650151851", "21.091174364770623"),
    @("586.151637066515", "22.814157353418754"),
    @("583.7078431656653", "25.663378578350205"),
    @("584.2085167228786", "25.26291670483937"),
    @("580.5063085183023", "34.37682038583849"),
    @("580.6443470713135", "43.831348506087735"),
    @("578.8400937301277", "38.34077282950777")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)

    # Force the numeric-looking literal to be stored as a text (shared-string)
    # cell rather than a number, without leaving any residual cell styling:
    # put it in as a formula returning a text literal, then convert the
    # formula result to a plain value in place via copy / paste-special.
    $bCell.Formula = "=""" + $data[$i][0] + """"
    $bCell.Copy()
    $bCell.PasteSpecial(-4163)

    $cCell.Formula = "=""" + $data[$i][1] + """"
    $cCell.Copy()
    $cCell.PasteSpecial(-4163)

    $ws.Cells.Item($row, 4).Value = "lorentzian"
}
